# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the per-language report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 00:50:22"
$wsZhCn.Range("H2").Value = "2016-03-19 00:50:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 00:50:26"
$wsDeDe.Range("H2").Value = "2016-03-19 00:50:46"
